$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17 currently store the Rotation value in column D (the "Layer"
# column). Move that value over to column E (Rotation) and replace
# column D with the literal layer name "Top".
for ($r = 2; $r -le 17; $r++) {
    $rotation = $ws.Cells.Item($r, 4).Value2
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.ClearFormats()
    $eCell.Value2 = $rotation
    $ws.Cells.Item($r, 4).Value = "Top"
}

# Update the active selection to match the author's saved cursor position.
$ws.Range("E20").Select()
